$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a brand-new "2022-Q3" worksheet right after "总计" (i.e. right
#    before "2022-Q2"), using the "2022-Q2" sheet as a structural template
#    (same headers / layout / styles), then overwrite it with the new
#    quarter's figures.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy($template, $null)
$newSheet = $wb.Worksheets.Item("2022-Q2 (2)")
$newSheet.Name = "2022-Q3"

# The fund code / name / header row are unchanged; only the figures differ.
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "0.58"
$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "95.06"
$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "3.73"
$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.0216"
$newSheet.Range("H2").Value = 8

# ---------------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: insert a new row right under the
#    header for the new "2022-Q3" quarter, pushing every existing quarter
#    row down by one, and renumber the running index in column A.
# ---------------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")
$totals.Rows.Item(2).Insert()

# Copy formatting down from the row that just got pushed to row 3 so the
# new row 2 matches the existing look (index-column style etc).
$totals.Range("A3:D3").Copy()
$totals.Range("A2:D2").PasteSpecial(-4122)  # xlPasteFormats

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q3"
$totals.Range("C2").Value = 1
$totals.Range("D2").Value = 0.02

for ($r = 3; $r -le 9; $r++) {
    $totals.Cells.Item($r, 1).Value = $r - 2
}
